$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: merge the new iteration-group header blocks first (so the
#     subsequent format copy is what finally sets their style, matching
#     the existing B1:D1 header block). ---
$ws.Range("E1:G1").Merge()
$ws.Range("H1:J1").Merge()

# --- Apply the existing header style (style already used by B1:D1 / B2:D2)
#     to the new header cells. Copy whole same-width row blocks so the
#     engine's border auto-merge logic doesn't mint extra style variants. ---
$ws.Range("B1:D1").Copy() | Out-Null
$ws.Range("E1:G1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("H1:J1").PasteSpecial(-4122) | Out-Null

$ws.Range("B2:D2").Copy() | Out-Null
$ws.Range("E2:G2").PasteSpecial(-4122) | Out-Null
$ws.Range("H2:J2").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Row 1: set the merged iteration-group header labels ---
$ws.Range("E1").Value = "Iteration_1"
$ws.Range("H1").Value = "Iteration_2"

# --- Row 2: repeat the 2030/2040/2050 interval labels for each group.
#     Copy the literal text values from B2/C2/D2 (xlPasteValues) instead of
#     assigning string literals, so "2030" stays text instead of being
#     auto-coerced to a number. ---
$ws.Range("B2").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4163) | Out-Null   # xlPasteAll
$ws.Range("H2").PasteSpecial(-4163) | Out-Null

$ws.Range("C2").Copy() | Out-Null
$ws.Range("F2").PasteSpecial(-4163) | Out-Null
$ws.Range("I2").PasteSpecial(-4163) | Out-Null

$ws.Range("D2").Copy() | Out-Null
$ws.Range("G2").PasteSpecial(-4163) | Out-Null
$ws.Range("J2").PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = $false

# --- Updated values in the pre-existing columns ---
$ws.Range("C4").Value = 2625.168436700786
$ws.Range("C5").Value = 1177664.151661842
$ws.Range("D7").Value = 1180474.661760693

# --- Row 4: Conventional ---
$ws.Range("E4").Value = 1183999.999999789
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 1184000.000000012
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0

# --- Row 5: Carbon Capture ---
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0

# --- Row 6: Electrification ---
$ws.Range("E6").Value = 0.0000006977175300198235
$ws.Range("F6").Value = 1181738.26257592
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0.0000008747791622340428
$ws.Range("I6").Value = 1181738.895756663
$ws.Range("J6").Value = 0

# --- Row 7: Water electrolysis ---
$ws.Range("E7").Value = 0.00000002113276659554965
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 1183903.89879729
$ws.Range("H7").Value = -0.00000002758638878031571
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 1183911.436787312

# --- Rows 8-15: all remaining new columns are zero ---
foreach ($r in 8..15) {
    foreach ($col in @("E","F","G","H","I","J")) {
        $ws.Range("$col$r").Value = 0
    }
}
